$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K10").Value = 5.689398
$ws.Range("L10").Value = 5.689398
$ws.Range("M10").Value = 5.689398
$ws.Range("N10").Value = 6.011680999999999
$ws.Range("O10").Value = 6.011680999999999
$ws.Range("P10").Value = 7.719391
$ws.Range("Q10").Value = 7.719391
$ws.Range("R10").Value = 8.082537
$ws.Range("S10").Value = 8.082537
